# Update "Fresh bloom Flowers_2025-9-25.xlsx"
# - Append 10 new flower rows (42-51) to the Orders sheet (column C)
# - Extend the dimension / ignoredErrors range on the Orders sheet to A1:L51
# - Extend the long numeric-as-text code in Summary!G2 with more trailing zeros

$wb = $excel.ActiveWorkbook

$wsOrders = $wb.Worksheets.Item("Orders")
$wsSummary = $wb.Worksheets.Item("Summary")

$newRows = @(
    "572_乒乓菊白_undefined_undefined_1bunch",
    "573_乒乓菊粉_undefined_undefined_1bunch",
    "734_乒乓菊红_undefined_undefined_1bunch",
    "15_深紫洋桔梗_Dark Purple Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g",
    "14_波浪浅紫洋桔梗_Wavy Light Purple Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g",
    "13_酒红洋桔梗_Burgundy Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g",
    "11_香槟洋桔梗_Champagne Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g",
    "12_肉粉洋桔梗_Peach Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g",
    "424_鼠尾白色_veronica white_undefined_1bunch",
    "423_鼠尾紫色_veronica purple_undefined_1bunch"
)

$startRow = 42
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $wsOrders.Cells.Item($row, 3).Value = $newRows[$i]
}

# Extend the long text code in Summary!G2 with 10 extra trailing zeros.
# Force a text number format first so the long numeric string isn't
# coerced into a floating point number (which would lose precision and
# the leading zero), then restore the default "Normal" style so no
# extra cell formatting is left behind.
$gCell = $wsSummary.Range("G2")
$gCell.NumberFormat = "@"
$gCell.Value = "0520000000000000000000000000000000000000000000000000"
$gCell.Style = "Normal"
